# Auto-generated update script applying the 2026-02-08 07:50 automatic data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/date/measurement values: safe to assign directly; Excel keeps them as text ---
$ws.Range('E2').Value = '2026-02-08 07:48:42'
$ws.Range('M2').Value = '-0.5 °C 7:23 TU'
$ws.Range('O2').Value = '-2.5 °C'
$ws.Range('E3').Value = '2026-02-08 07:48:44'
$ws.Range('E4').Value = '2026-02-08 07:48:47'
$ws.Range('J4').Value = '1001.4 hPa'
$ws.Range('K4').Value = '0.0 MJ/m2'
$ws.Range('N4').Value = '4.4 °C 7:17 TU'
$ws.Range('O4').Value = '7.7 °C'
$ws.Range('E5').Value = '2026-02-08 07:48:50'
$ws.Range('E6').Value = '2026-02-08 07:48:52'
$ws.Range('J6').Value = '1001.2 hPa'
$ws.Range('O6').Value = '7.9 °C'
$ws.Range('E7').Value = '2026-02-08 07:48:55'
$ws.Range('N7').Value = '9.9 °C 7:15 TU'
$ws.Range('O7').Value = '10.9 °C'
$ws.Range('E8').Value = '2026-02-08 07:48:57'
$ws.Range('J8').Value = '1001.4 hPa'
$ws.Range('E9').Value = '2026-02-08 07:49:00'
$ws.Range('E10').Value = '2026-02-08 07:49:02'
$ws.Range('N10').Value = '3.3 °C 7:09 TU'
$ws.Range('O10').Value = '6.4 °C'
$ws.Range('E11').Value = '2026-02-08 07:49:04'
$ws.Range('N11').Value = '-0.8 °C 7:28 TU'
$ws.Range('O11').Value = '1.2 °C'
$ws.Range('E12').Value = '2026-02-08 07:49:07'
$ws.Range('N12').Value = '6.7 °C 7:03 TU'
$ws.Range('O12').Value = '8.4 °C'
$ws.Range('E13').Value = '2026-02-08 07:49:10'
$ws.Range('N13').Value = '-1.5 °C 7:09 TU'
$ws.Range('O13').Value = '0.2 °C'
$ws.Range('E14').Value = '2026-02-08 07:49:12'
$ws.Range('N14').Value = '5.8 °C 7:28 TU'
$ws.Range('O14').Value = '7.9 °C'
$ws.Range('E15').Value = '2026-02-08 07:49:15'
$ws.Range('O15').Value = '6.4 °C'
$ws.Range('E16').Value = '2026-02-08 07:49:17'
$ws.Range('L16').Value = '34.6 km/h - 166º 7:27 TU'
$ws.Range('E17').Value = '2026-02-08 07:49:20'
$ws.Range('G17').Value = '2 cm'
$ws.Range('K17').Value = '0.1 MJ/m2'
$ws.Range('N17').Value = '-1.5 °C 7:03 TU'
$ws.Range('E18').Value = '2026-02-08 07:49:22'
$ws.Range('N18').Value = '4.2 °C 7:22 TU'
$ws.Range('O18').Value = '7.3 °C'
$ws.Range('E19').Value = '2026-02-08 07:49:24'
$ws.Range('E20').Value = '2026-02-08 07:49:27'
$ws.Range('K20').Value = '0.1 MJ/m2'
$ws.Range('E21').Value = '2026-02-08 07:49:29'
$ws.Range('L21').Value = '8.3 km/h - 197º 7:24 TU'
$ws.Range('N21').Value = '1.6 °C 7:21 TU'
$ws.Range('E22').Value = '2026-02-08 07:49:31'
$ws.Range('N22').Value = '-8.3 °C 7:01 TU'
$ws.Range('O22').Value = '-6.9 °C'
$ws.Range('E23').Value = '2026-02-08 07:49:34'
$ws.Range('E24').Value = '2026-02-08 07:49:37'
$ws.Range('N24').Value = '2.3 °C 7:09 TU'
$ws.Range('O24').Value = '5.9 °C'
$ws.Range('E25').Value = '2026-02-08 07:49:39'
$ws.Range('L25').Value = '25.6 km/h - 345º 7:22 TU'
$ws.Range('E26').Value = '2026-02-08 07:49:42'
$ws.Range('N26').Value = '-0.2 °C 7:01 TU'
$ws.Range('E27').Value = '2026-02-08 07:49:45'
$ws.Range('O27').Value = '-4.0 °C'
$ws.Range('E28').Value = '2026-02-08 07:49:47'
$ws.Range('J28').Value = '1001.6 hPa'
$ws.Range('N28').Value = '2.3 °C 7:22 TU'
$ws.Range('O28').Value = '5.2 °C'
$ws.Range('E29').Value = '2026-02-08 07:49:50'
$ws.Range('O29').Value = '8.7 °C'
$ws.Range('E30').Value = '2026-02-08 07:49:53'
$ws.Range('N30').Value = '6.6 °C 7:19 TU'
$ws.Range('O30').Value = '8.8 °C'
$ws.Range('E31').Value = '2026-02-08 07:49:55'
$ws.Range('J31').Value = '999.7 hPa'
$ws.Range('E32').Value = '2026-02-08 07:49:58'
$ws.Range('L32').Value = '11.9 km/h - 300º 7:20 TU'
$ws.Range('O32').Value = '2.0 °C'
$ws.Range('E33').Value = '2026-02-08 07:50:01'
$ws.Range('J33').Value = '1003.5 hPa'
$ws.Range('N33').Value = '-1.4 °C 7:02 TU'
$ws.Range('O33').Value = '0.2 °C'
$ws.Range('E34').Value = '2026-02-08 07:50:04'
$ws.Range('O34').Value = '-1.2 °C'
$ws.Range('E35').Value = '2026-02-08 07:50:06'
$ws.Range('K35').Value = '0.0 MJ/m2'
$ws.Range('E36').Value = '2026-02-08 07:50:09'
$ws.Range('N36').Value = '8.4 °C 7:10 TU'
$ws.Range('O36').Value = '10.4 °C'
$ws.Range('E37').Value = '2026-02-08 07:50:12'
$ws.Range('N37').Value = '0.7 °C 7:29 TU'
$ws.Range('O37').Value = '2.9 °C'
$ws.Range('E38').Value = '2026-02-08 07:50:14'
$ws.Range('K38').Value = '0.0 MJ/m2'
$ws.Range('N38').Value = '4.0 °C 7:02 TU'
$ws.Range('O38').Value = '6.9 °C'
$ws.Range('E39').Value = '2026-02-08 07:50:17'
$ws.Range('E40').Value = '2026-02-08 07:50:20'
$ws.Range('J40').Value = '1004.2 hPa'
$ws.Range('N40').Value = '0.8 °C 7:08 TU'
$ws.Range('O40').Value = '2.4 °C'
$ws.Range('E41').Value = '2026-02-08 07:50:22'
$ws.Range('E42').Value = '2026-02-08 07:50:25'
$ws.Range('O42').Value = '8.7 °C'
$ws.Range('E43').Value = '2026-02-08 07:50:27'
$ws.Range('O43').Value = '4.0 °C'
$ws.Range('E44').Value = '2026-02-08 07:50:29'
$ws.Range('E45').Value = '2026-02-08 07:50:32'
$ws.Range('J45').Value = '1002.3 hPa'
$ws.Range('K45').Value = '0.0 MJ/m2'
$ws.Range('O45').Value = '2.4 °C'
$ws.Range('E46').Value = '2026-02-08 07:50:35'
$ws.Range('J46').Value = '1002.1 hPa'
$ws.Range('K46').Value = '0.0 MJ/m2'
$ws.Range('O46').Value = '5.9 °C'

# --- Percentage-looking values ("NN%"): Excel auto-converts these to numbers with a percent
# number format when assigned directly. Stage each one through a scratch cell formatted as
# Text ("@"), copy/paste-special the resulting text value into the target so the target
# keeps its original style/number format, then fully clear the scratch cell so it does not
# linger in the worksheets used range.
$scratch = $ws.Range("ZZ9000")
$percentUpdates = [ordered]@{
  'H2' = '89%'
  'H4' = '75%'
  'H5' = '83%'
  'H9' = '78%'
  'H12' = '79%'
  'H18' = '81%'
  'H28' = '83%'
  'H30' = '69%'
  'H31' = '63%'
  'H32' = '99%'
  'H36' = '74%'
  'H37' = '89%'
  'H38' = '86%'
  'H39' = '86%'
  'H42' = '91%'
  'H46' = '89%'
}
foreach ($cellRef in $percentUpdates.Keys) {
  $scratch.NumberFormat = "@"
  $scratch.Value = $percentUpdates[$cellRef]
  $scratch.Copy()
  $target = $ws.Range($cellRef)
  $target.PasteSpecial(-4163)
  $scratch.Clear()
}

Write-Host "Applied automatic data/banner refresh for 2026-02-08 07:50"
